# Rewrite the 시간대 (hour-of-day) counts table: re-sort the 24 hour
# buckets into the new descending-by-count order and apply the new
# counts.
#
# Row 1 (header) is untouched: A1 = 접수시간대, B1 = 0.
# Rows 2-25 get the new hour label (text, zero-padded) in column A and
# the new count (number) in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, new hour label, whether the label differs from what's already there, new count
$rows = @(
    @(2,  "08", $true,  70),
    @(3,  "11", $true,  61),
    @(4,  "09", $true,  58),
    @(5,  "10", $true,  51),
    @(6,  "24", $true,  44),
    @(7,  "13", $true,  44),
    @(8,  "14", $true,  32),
    @(9,  "07", $true,  31),
    @(10, "16", $true,  30),
    @(11, "17", $true,  28),
    @(12, "15", $false, 23),
    @(13, "18", $false, 16),
    @(14, "21", $true,  16),
    @(15, "20", $true,  14),
    @(16, "06", $true,  13),
    @(17, "22", $false, 12),
    @(18, "05", $true,  11),
    @(19, "23", $true,  11),
    @(20, "19", $true,  11),
    @(21, "12", $true,  9),
    @(22, "01", $true,  9),
    @(23, "02", $true,  7),
    @(24, "03", $true,  5),
    @(25, "04", $true,  2)
)

foreach ($row in $rows) {
    $r = $row[0]
    $label = $row[1]
    $labelChanged = $row[2]
    $count = $row[3]

    if ($labelChanged) {
        # Format as Text first so the zero-padded hour label ("08", "09",
        # ...) is preserved as a string instead of being auto-converted
        # to a number.
        $ws.Cells.Item($r, 1).NumberFormat = "@"
        $ws.Cells.Item($r, 1).Value = $label
    }

    $ws.Cells.Item($r, 2).Value = $count
}
